# Generate Report for Handoff
# A new handoff was generated for "a9fba91c-9be6-4a20-b0e5-dbd4fc4697cd.md"
# (the 3rd data row on every sheet). Update its Status to "Ready for handoff"
# and refresh the Latest Handoff Datetime values accordingly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (summary) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 22:53:32"

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 22:53:28"

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 22:53:32"
